$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "Förändrad" (changed) date column C for rows 2-8 from
# 2023-09-05 (serial 45174) to 2023-09-06 (serial 45175), preserving
# the existing cell formatting/style.
foreach ($r in 2..8) {
    $ws.Cells.Item($r, 3).Value = 45175
}
